$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from column F (old column D) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Populate the two new quarterly-data columns (D, E) with their values
$ws.Range("D7").Value = 43460
$ws.Range("E7").Value = 43369
$ws.Range("D8").Value = 790700
$ws.Range("E8").Value = 753800
$ws.Range("D9").Value = 667400
$ws.Range("E9").Value = 647200
$ws.Range("D10").Value = 123300
$ws.Range("E10").Value = 106600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3900
$ws.Range("E14").Value = 3300
$ws.Range("D15").Value = 36100
$ws.Range("E15").Value = 37000
$ws.Range("D17").Value = 741100
$ws.Range("E17").Value = 706900
$ws.Range("D18").Value = 49600
$ws.Range("E18").Value = 46900
$ws.Range("D20").Value = 800
$ws.Range("E20").Value = 800
$ws.Range("D21").Value = 86500
$ws.Range("E21").Value = 84700
$ws.Range("D22").Value = 15400
$ws.Range("E22").Value = 15600
$ws.Range("D23").Value = 35000
$ws.Range("E23").Value = 32100
$ws.Range("D24").Value = 3000
$ws.Range("E24").Value = 5700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 32000
$ws.Range("E26").Value = 26400
$ws.Range("D27").Value = 32000
$ws.Range("E27").Value = 26400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("E32").Value = -800
$ws.Range("D33").Value = 32000
$ws.Range("E33").Value = 26400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 32000
$ws.Range("E35").Value = 26400
$ws.Range("D38").Value = 43460
$ws.Range("E38").Value = 43369
$ws.Range("D41").Value = 16200
$ws.Range("E41").Value = 11000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 85900
$ws.Range("E43").Value = 42700
$ws.Range("D44").Value = 24000
$ws.Range("E44").Value = 23400
$ws.Range("D45").Value = 69900
$ws.Range("E45").Value = 73900
$ws.Range("D46").Value = 196000
$ws.Range("E46").Value = 151000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 769300
$ws.Range("E48").Value = 762200
$ws.Range("D49").Value = 186700
$ws.Range("E49").Value = 187600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 142800
$ws.Range("E52").Value = 143200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1294800
$ws.Range("E54").Value = 1244000
$ws.Range("D57").Value = 112200
$ws.Range("E57").Value = 97200
$ws.Range("D58").Value = 8100
$ws.Range("E58").Value = 7400
$ws.Range("D59").Value = 367700
$ws.Range("E59").Value = 403000
$ws.Range("D60").Value = 488000
$ws.Range("E60").Value = 507600
$ws.Range("D61").Value = 1263900
$ws.Range("E61").Value = 1153000
$ws.Range("D62").Value = 398100
$ws.Range("E62").Value = 399300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2150000
$ws.Range("E66").Value = 2059900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2704000
$ws.Range("E72").Value = 2686500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -855200
$ws.Range("E76").Value = -815900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43460
$ws.Range("E80").Value = 43369
$ws.Range("D81").Value = 32000
$ws.Range("E81").Value = 26400
$ws.Range("D83").Value = 36100
$ws.Range("E83").Value = 37000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 6600
$ws.Range("E89").Value = 49600
$ws.Range("D91").Value = -47500
$ws.Range("E91").Value = -31200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -35300
$ws.Range("E94").Value = 418500
$ws.Range("D96").Value = -15400
$ws.Range("E96").Value = -16200
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 33900
$ws.Range("E100").Value = -468000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 5200
$ws.Range("E102").Value = 100

# A handful of historical quarters were also restated/corrected during this update
$ws.Range("F17").Value = 746700
$ws.Range("F18").Value = 70400
$ws.Range("F20").Value = 800
$ws.Range("F32").Value = -800
$ws.Range("J91").Value = -22800
